$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-18 Thursday" "2024-04-19 Friday"

Replace-Text "152÷4=" "669÷3="
Replace-Text "153÷3=" "918÷4="
Replace-Text "945÷8=" "734÷9="
Replace-Text "416÷9=" "108÷7="
Replace-Text "857÷6=" "757÷5="

Replace-Text "566÷6=" "124÷8="
Replace-Text "409÷6=" "535÷3="
Replace-Text "732÷8=" "675÷6="
Replace-Text "503÷2=" "596÷3="
Replace-Text "475÷3=" "775÷8="

Replace-Text "156÷5=" "835÷9="
Replace-Text "838÷3=" "874÷2="
Replace-Text "911÷8=" "728÷8="
Replace-Text "916÷2=" "856÷4="
Replace-Text "953÷2=" "389÷9="

Replace-Text "749÷2=" "162÷7="
Replace-Text "910÷8=" "956÷3="
Replace-Text "330÷2=" "868÷3="
Replace-Text "862÷9=" "881÷9="
Replace-Text "507÷8=" "467÷7="

Replace-Text "297÷3=" "852÷2="
Replace-Text "746÷4=" "337÷7="
Replace-Text "631÷8=" "963÷2="
Replace-Text "395÷6=" "154÷2="
Replace-Text "479÷3=" "926÷9="
